$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = -7
$ws.Range("F5").Value = 5
$ws.Range("F6").Value = 8
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = -2
